$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.813.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.89%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.942.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.75%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'242.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.95%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  -0.08%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4890"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.28%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.2958"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.27%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.06881"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.63%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'19.41"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.64%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'106.27"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'1.945.19"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.55%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.07719"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.37%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'5.352"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.20%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.6983"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -2.01%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'277.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.65%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'30.817.12"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.70%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.000007712"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.57%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'13.13"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.74%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "'2.196.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.19%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("B21").Value = "'Dai"
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").Value = "'1.000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.02%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'5.470"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.32%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.9998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.05%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'6.536"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.77%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -2.08%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'167.62"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.36%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'19.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.44%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -0.91%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.1048"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.86%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'1.391"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -3.46%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("B31").Value = "'Filecoin"
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = "'4.565"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -3.99%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("B32").Value = "'PancakeSwap"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "'1.554"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -2.76%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'4.372"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -3.06%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.04856"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -2.85%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.7535"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.75%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -0.62%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -0.01%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'2.732"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.09%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.01996"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -2.66%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.653"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.93%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'6.526"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +1.25%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'77.97"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +6.81%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -2.18%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.9079"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +2.81%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'108.23"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.34%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.4397"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.75%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.9989"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.17%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'7.764"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +3.43%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'998.04"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.22%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -1.89%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'9.252"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.84%  "
$ws.Range("E51").Style = "Normal"
